$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

$ws.Range("D5").Value = "KONTOSTAND AM 30.06.2024"

$ws.Range("B6").Value = "03.07."
$ws.Range("C6").Value = "04.07."
$ws.Range("D6").Value = "BEITRAG Allianz SE K-10722420"
$ws.Range("E6").Value = "55,07-"

$ws.Range("B7").Value = "04.07."
$ws.Range("C7").Value = "05.07."
$ws.Range("D7").Value = "KARTENZ./04.07 ALDI SUED RO"
$ws.Range("E7").Value = "81,56-"

$ws.Range("B8").Value = "08.07."
$ws.Range("C8").Value = "09.07."
$ws.Range("D8").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E8").Value = "24,50-"

$ws.Range("B9").Value = "12.07."
$ws.Range("C9").Value = "13.07."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 69717726"
$ws.Range("E9").Value = "85,19-"

$ws.Range("B10").Value = "13.07."
$ws.Range("C10").Value = "14.07."
$ws.Range("D10").Value = "RECHNUNG VODAFONE GMBH 46256545"
$ws.Range("E10").Value = "40,24-"

$ws.Range("D12").Value = "KONTOSTAND AM 17.07.2024"
$ws.Range("E12").Value = "286,56-"

$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 24.07.2024"
